$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full 3x3 cross of sending/target clusters (ECs, FAPs, sCs) for Mdk-Tspan1 pair
$data = @(
  @(2, "ECs", "ECs", 3, 1, 2.891504666666667, 8.674514, 0.1213590456377548, 0.1213590456377548, 1, 0.3333333333333333, 0.18316, 0.54948, 0.07387220463520254, 0.07387220463520254, 0.5296079947466666, 4.76647195272, 0.008965060253685105, 0.008965060253685105),
  @(3, "ECs", "FAPs", 3, 1, 2.891504666666667, 8.674514, 0.1213590456377548, 0.1213590456377548, 3, 1, 1.360168666666667, 4.080506, 0.5485840690237529, 0.548584069023753, 3.932934047120444, 35.396406424084, 0.06657563906879883, 0.06657563906879885),
  @(4, "ECs", "sCs", 3, 1, 2.891504666666667, 8.674514, 0.1213590456377548, 0.1213590456377548, 3, 1, 0.9360883333333335, 2.808265, 0.3775437263410444, 0.3775437263410445, 2.706703784245556, 24.36033405821, 0.04581834631527081, 0.04581834631527081),
  @(5, "FAPs", "ECs", 3, 1, 12.04042966666667, 36.121289, 0.505347637947847, 0.505347637947847, 1, 0.3333333333333333, 0.18316, 0.54948, 0.07387220463520254, 0.07387220463520254, 2.205325097746667, 19.84792587972, 0.03733114412239959, 0.03733114412239959),
  @(6, "FAPs", "FAPs", 3, 1, 12.04042966666667, 36.121289, 0.505347637947847, 0.505347637947847, 3, 1, 1.360168666666667, 4.080506, 0.5485840690237529, 0.548584069023753, 16.37701516580378, 147.393136492234, 0.2772256634969722, 0.2772256634969723),
  @(7, "FAPs", "sCs", 3, 1, 12.04042966666667, 36.121289, 0.505347637947847, 0.505347637947847, 3, 1, 0.9360883333333335, 2.808265, 0.3775437263410444, 0.3775437263410445, 11.27090573928723, 101.438151653585, 0.1907908303284752, 0.1907908303284752),
  @(8, "sCs", "ECs", 3, 1, 8.894099000000001, 26.682297, 0.3732933164143983, 0.3732933164143982, 1, 0.3333333333333333, 0.18316, 0.54948, 0.07387220463520254, 0.07387220463520254, 1.62904317284, 14.66138855556, 0.02757600025911784, 0.02757600025911784),
  @(9, "sCs", "FAPs", 3, 1, 8.894099000000001, 26.682297, 0.3732933164143983, 0.3732933164143982, 3, 1, 1.360168666666667, 4.080506, 0.5485840690237529, 0.548584069023753, 12.09747477803133, 108.877273002282, 0.2047827664579819, 0.2047827664579819),
  @(10, "sCs", "sCs", 3, 1, 8.894099000000001, 26.682297, 0.3732933164143983, 0.3732933164143982, 3, 1, 0.9360883333333335, 2.808265, 0.3775437263410444, 0.3775437263410445, 8.325662309411669, 74.93096078470502, 0.1409345496972985, 0.1409345496972985)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]   # A: Sending cluster
    $ws.Cells.Item($r, 2).Value = "Mdk"     # B: Ligand symbol
    $ws.Cells.Item($r, 3).Value = "Tspan1"  # C: Receptor symbol
    $ws.Cells.Item($r, 4).Value = $row[2]   # D: Target cluster
    for ($i = 0; $i -lt 16; $i++) {
        $ws.Cells.Item($r, 5 + $i).Value = $row[3 + $i]
    }
}